$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.855.21'
$ws.Range("E2").Value = '  -0.16%  '

$ws.Range("D3").Value = '1.870.01'
$ws.Range("E3").Value = '  -1.30%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9987'
$ws.Range("E4").Value = '  -0.24%  '

$ws.Range("E5").Value = '  -4.00%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '241.83'
$ws.Range("E6").Value = '  -1.18%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9979'
$ws.Range("E7").Value = '  -0.32%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3151'
$ws.Range("E8").Value = '  +0.59%  '

$ws.Range("B9").Value = 'Solana'
$ws.Range("C9").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '24.70'
$ws.Range("E9").Value = '  -3.79%  '

$ws.Range("B10").Value = 'Dogecoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07126'
$ws.Range("E10").Value = '  -1.80%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08408'
$ws.Range("E11").Value = '  -5.47%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.7532'
$ws.Range("E12").Value = '  -2.49%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.429'
$ws.Range("E13").Value = '  -0.26%  '

$ws.Range("D14").Value = '1.835.70'
$ws.Range("E14").Value = '  -1.85%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '92.67'
$ws.Range("E15").Value = '  -1.80%  '

$ws.Range("D16").Value = '29.842.79'
$ws.Range("E16").Value = '  -0.06%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.036'
$ws.Range("E17").Value = '  -2.53%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.58'
$ws.Range("E18").Value = '  -2.79%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '243.10'
$ws.Range("E19").Value = '  -1.12%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007812'
$ws.Range("E20").Value = '  -0.65%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9984'
$ws.Range("E21").Value = '  -0.20%  '

$ws.Range("D22").Value = '2.116.47'
$ws.Range("E22").Value = '  +0.70%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.955'
$ws.Range("E23").Value = '  -2.20%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.9957'
$ws.Range("E24").Value = '  -0.54%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1583'
$ws.Range("E25").Value = '  -1.01%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.301'
$ws.Range("E26").Value = '  -2.22%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '164.08'
$ws.Range("E27").Value = '  +0.67%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.58'
$ws.Range("E28").Value = '  -1.25%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.025'
$ws.Range("E29").Value = '  -1.07%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.471'
$ws.Range("E30").Value = '  +3.07%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.649'
$ws.Range("E31").Value = '  +2.22%  '

$ws.Range("E32").Value = '  -0.99%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.297'
$ws.Range("E33").Value = '  +4.67%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05324'
$ws.Range("E34").Value = '  -3.41%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.236'
$ws.Range("E35").Value = '  -1.00%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7528'
$ws.Range("E36").Value = '  +0.10%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.000'
$ws.Range("E37").Value = '  +0.51%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.696'
$ws.Range("E38").Value = '  -0.83%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.745'
$ws.Range("E40").Value = '  -1.67%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.4471'
$ws.Range("E41").Value = '  -0.69%  '

$ws.Range("D42").Value = '1.109.81'
$ws.Range("E42").Value = '  +2.09%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.087'
$ws.Range("E43").Value = '  +0.69%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '72.19'
$ws.Range("E44").Value = '  -2.04%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8598'
$ws.Range("E45").Value = '  +0.53%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.9991'
$ws.Range("E46").Value = '  -0.16%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '103.24'
$ws.Range("E47").Value = '  +0.59%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.690'
$ws.Range("E48").Value = '  +1.26%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.844'
$ws.Range("E49").Value = '  -2.37%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.045'
$ws.Range("E50").Value = '  +1.81%  '

$ws.Range("D51").Value = '2.012.80'
$ws.Range("E51").Value = '  +0.81%  '
